# wspr_calc_direct_shift_lower_pll_range.xlsx
#
# Commit: "optimal denom for pll 500 (spreadsheet says good for 10/12/15/17/20M)"
#
# Before: 3 sheets  -> "pll 390", "pll 600", "pll 700"
# After : 5 sheets  -> "pll 500", "pll 416", "pll 390", "pll 600", "pll 700"
#
# The two new sheets are duplicates of the existing "pll 390" and "pll 600"
# sheets (same layout/formulas/hyperlink), inserted in front of "pll 390",
# then re-tuned with new PLL divisor values in row 7 / row 8.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate "pll 390" -> becomes "pll 500" (placed right before "pll 390")
# ---------------------------------------------------------------------
$pll390 = $wb.Worksheets.Item("pll 390")
$pll390.Copy($pll390)
$pll500 = $wb.Worksheets.Item("pll 390 (2)")
$pll500.Name = "pll 500"

# ---------------------------------------------------------------------
# 2) Duplicate "pll 600" -> becomes "pll 416" (placed right before "pll 390",
#    i.e. right after the newly created "pll 500")
# ---------------------------------------------------------------------
$pll600 = $wb.Worksheets.Item("pll 600")
$pll390again = $wb.Worksheets.Item("pll 390")
$pll600.Copy($pll390again)
$pll416 = $wb.Worksheets.Item("pll 600 (2)")
$pll416.Name = "pll 416"

# Final tab order is now: pll 500, pll 416, pll 390, pll 600, pll 700

# ---------------------------------------------------------------------
# 3) Re-tune "pll 500" divisors (row 7 = 2nd-stage divisor candidates,
#    row 8 = chosen divisor per column-pair)
# ---------------------------------------------------------------------
$pll500 = $wb.Worksheets.Item("pll 500")
$pll500.Range("D7").Value2 = 18
$pll500.Range("F7").Value2 = 20
$pll500.Range("H7").Value2 = 24
$pll500.Range("J7").Value2 = 28
$pll500.Range("L7").Value2 = 36

$pll500.Range("D8").Value2 = 19
$pll500.Range("F8").Value2 = 19
$pll500.Range("H8").Value2 = 19
$pll500.Range("J8").Value2 = 19
$pll500.Range("L8").Value2 = 19

# ---------------------------------------------------------------------
# 4) Re-tune "pll 416" divisors
# ---------------------------------------------------------------------
$pll416 = $wb.Worksheets.Item("pll 416")
$pll416.Range("D7").Value2 = 14
$pll416.Range("F7").Value2 = 16
$pll416.Range("H7").Value2 = 20
$pll416.Range("J7").Value2 = 22
$pll416.Range("L7").Value2 = 30

$pll416.Range("D8").Value2 = 15
$pll416.Range("F8").Value2 = 15
$pll416.Range("H8").Value2 = 16
$pll416.Range("J8").Value2 = 15
$pll416.Range("L8").Value2 = 16

# ---------------------------------------------------------------------
# 5) Selections per sheet + which tab ends up active ("pll 600", index 3)
# ---------------------------------------------------------------------
$pll500 = $wb.Worksheets.Item("pll 500")
$pll500.Activate()
$pll500.Range("C19").Select()

$pll416 = $wb.Worksheets.Item("pll 416")
$pll416.Activate()
$pll416.Range("A14").Select()

$pll390 = $wb.Worksheets.Item("pll 390")
$pll390.Activate()
$pll390.Range("A13").Select()

$pll700 = $wb.Worksheets.Item("pll 700")
$pll700.Activate()
$pll700.Range("E9").Select()

$pll600 = $wb.Worksheets.Item("pll 600")
$pll600.Activate()
$pll600.Range("J13").Select()
